$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 43 is a new data row appended after the existing A1:D42 data block.
# Column A holds a date-formatted string ("2025/10/01"). Assigning that
# literal via .Value would be auto-coerced into a date serial by Excel's
# normal type inference, so the cell is pre-formatted as Text, the literal
# text is written, and the formatting is cleared again afterwards so the
# cell ends up with no explicit style - identical to how the sibling data
# cells in the rest of the table are stored (plain text, default style).
$ws.Range("A43").NumberFormat = "@"
$ws.Range("A43").Value = "2025/10/01"
$ws.Range("A43").ClearFormats()

$ws.Range("B43").Value = "水"
$ws.Range("C43").Value = 6
$ws.Range("D43").Value = 3
